$wb = $excel.ActiveWorkbook
$wsOld = $wb.Worksheets.Item("Previously added")
$wsNew = $wb.Worksheets.Item("New")

# ---------------------------------------------------------------------------
# 1) "Previously added" sheet: append the record that currently sits in row 2
#    of the "New" sheet as the new last row (354), with a matching hyperlink.
# ---------------------------------------------------------------------------

$oldLink = "https://www.ss.com/msg/lv/real-estate/wood/balvi-and-reg/berzpils-pag/mnmcj.html"

$wsOld.Range("A354").Value = $oldLink
$wsOld.Hyperlinks.Add($wsOld.Range("A354"), $oldLink)
$wsOld.Range("B354").Value = "11 000 €"
$wsOld.Range("C354").Value = "Balvi un raj."
$wsOld.Range("D354").Value = "2 ha."
$wsOld.Range("E354").NumberFormat = "@"
$wsOld.Range("E354").Value = "38500050160"
$wsOld.Range("F354").Value = 46024.558333333334

# Re-apply the row formatting used by the rest of the table (copies the
# cell styles only, so the values set above are preserved) - this keeps
# the new row visually/structurally identical to its neighbours.
$wsOld.Range("A353:F353").Copy()
$wsOld.Range("A354:F354").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) "New" sheet: the old row 2 has just been relocated into "Previously
#    added", so wipe the hyperlink bookkeeping on this sheet and lay down
#    five freshly scraped rows (2-6).
# ---------------------------------------------------------------------------

# Stash the pristine row-2 cell formatting (link/text/date styles) in a
# scratch row before anything on row 2 gets touched - Hyperlinks.Add() and
# NumberFormat edits below would otherwise clobber the style indices.
$wsNew.Range("A2:F2").Copy()
$wsNew.Range("A200:F200").PasteSpecial(-4122)

$wsNew.Range("A200:F200").Copy()
$wsNew.Range("A3:F6").PasteSpecial(-4122)
$wsNew.Hyperlinks.Delete()

$rows = @(
    @{ row=2; link="https://www.ss.com/msg/lv/real-estate/wood/bauska-and-reg/islices-pag/ofhki.html"; price="6 300 €"; district="Bauska un raj."; area="2 ha."; cadastre=""; date=46025.58194444445 },
    @{ row=3; link="https://www.ss.com/msg/lv/real-estate/wood/jekabpils-and-reg/krustpils-pag/jmnen.html"; price="20 000 €"; district="Jēkabpils un raj."; area="3.65 ha."; cadastre="56680060283"; date=46025.65763888889 },
    @{ row=4; link="https://www.ss.com/msg/lv/real-estate/wood/jekabpils-and-reg/krustpils-pag/kfdjl.html"; price="40 000 €"; district="Jēkabpils un raj."; area="7 ha."; cadastre="56680060573"; date=46025.64722222222 },
    @{ row=5; link="https://www.ss.com/msg/lv/real-estate/wood/rezekne-and-reg/audrinu-pag/cxdpb.html"; price="6 200 €"; district="Rēzekne un raj."; area="1.10 ha."; cadastre="78420020154"; date=46026.76736111111 },
    @{ row=6; link="https://www.ss.com/msg/lv/real-estate/wood/rezekne-and-reg/lendzu-pag/blmofi.html"; price="75 000 €"; district="Rēzekne un raj."; area="3 ha."; cadastre="78660040107"; date=46026.73541666666 }
)

foreach ($r in $rows) {
    $rowNum = $r.row
    $wsNew.Range("A$rowNum").Value = $r.link
    $wsNew.Hyperlinks.Add($wsNew.Range("A$rowNum"), $r.link)
    $wsNew.Range("B$rowNum").Value = $r.price
    $wsNew.Range("C$rowNum").Value = $r.district
    $wsNew.Range("D$rowNum").Value = $r.area
    $wsNew.Range("E$rowNum").NumberFormat = "@"
    $wsNew.Range("E$rowNum").Value = $r.cadastre
    $wsNew.Range("F$rowNum").Value = $r.date
}

# Re-apply the original row-2 style (link/text/date formats) across all five
# rows so the newly added cells reuse the same style indices as before,
# instead of the ad-hoc "Hyperlink" style Hyperlinks.Add introduces.
$wsNew.Range("A200:F200").Copy()
$wsNew.Range("A2:F6").PasteSpecial(-4122)

# Drop the scratch row used to stash the formatting.
$wsNew.Range("A200:F200").Clear()
